# Add Asatru wonder Freya Sanctuary (AB4-350)
# Updates the "Religious Buildings" sheet Wonders summary table (rows 7, 11, 18)
# to reflect the new wonder totals, and moves the active selection to T7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Religious Buildings")

# Row 7 - Asatru: add new wonder stats (Curr/Rel.Happ./Obs.Comm./Rel.Ch./Bonus)
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 3
$ws.Range("T7").Value = 2

# Row 11 - Hellenism: updated aggregate wonder stats
$ws.Range("N11").Value = 2
$ws.Range("P11").Value = 6
$ws.Range("T11").Value = 3

# Row 18 - Zoroastrianism: updated aggregate wonder stats
$ws.Range("N18").Value = 3
$ws.Range("O18").Value = 6
$ws.Range("P18").Value = 14
$ws.Range("Q18").Value = 12
$ws.Range("T18").Value = 4

# Move the selection to T7, matching the saved view state
$ws.Activate()
$ws.Range("T7").Select()
